# Generate Report for Handoff
# - Updates the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview!E2:F2, zh-cn!C2, de-de!C2)
# - Refreshes the associated handoff timestamps
# - Narrows the now-shorter "Status" columns (they were previously auto-fit to the
#   longer text) on Overview (E:F) and on the zh-cn/de-de sheets (C)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update "Status" value wherever it is shown ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refresh the handoff generation timestamps ---
$wsOverview.Range("G2").Value = "2016-08-12 03:20:26"
$wsZhCn.Range("H2").Value     = "2016-08-12 03:20:21"
$wsDeDe.Range("H2").Value     = "2016-08-12 03:20:26"

# --- Shrink the "Status" columns to fit the new, shorter text ---
$wsOverview.Range("E1").ColumnWidth = 16.35
$wsOverview.Range("F1").ColumnWidth = 16.35
$wsZhCn.Range("C1").ColumnWidth     = 16.35
$wsDeDe.Range("C1").ColumnWidth     = 16.35
